# Applies the diff: swaps F:V data between pairs of adjacent match rows,
# and appends two new match rows (165, 166) at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap columns F:V between pairs of rows (same data_partida, rows got
#        reordered by the upstream scraper) ---
$pairs = @(
    @(22, 23),
    @(24, 25),
    @(47, 48),
    @(51, 52),
    @(101, 102),
    @(127, 128),
    @(129, 130),
    @(139, 140),
    @(142, 143),
    @(162, 163)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range($ws.Cells.Item($r1, 6), $ws.Cells.Item($r1, 22))
    $rng2 = $ws.Range($ws.Cells.Item($r2, 6), $ws.Cells.Item($r2, 22))

    $v1 = $rng1.Value2
    $v2 = $rng2.Value2

    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}

# --- 2. Append two new rows (165 and 166) with new match data ---

# Copy formatting (styles/number formats) from the last existing row (164)
# so the new rows keep the same look (bold/border index column, date format, etc.)
# Only copy the used A:V range (not the whole row) to avoid touching formatting
# of unused columns beyond V.
$srcRange = $ws.Range($ws.Cells.Item(164, 1), $ws.Cells.Item(164, 22))
$srcRange.Copy() | Out-Null
$dst165 = $ws.Range($ws.Cells.Item(165, 1), $ws.Cells.Item(165, 22))
$dst165.PasteSpecial(-4122) | Out-Null # xlPasteFormats
$dst166 = $ws.Range($ws.Cells.Item(166, 1), $ws.Cells.Item(166, 22))
$dst166.PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = 0

function Set-MatchRow {
    param(
        $RowNum,
        $Indice,
        $Home,
        $HomeGoals,
        $Away,
        $AwayGoals,
        $DataPartida,
        $HomeOpeningOdds,
        $HomeOpeningDataHora,
        $HomeClosingOdds,
        $HomeClosingDataHora,
        $DrawOpeningOdds,
        $DrawOpeningDataHora,
        $DrawClosingOdds,
        $DrawClosingDataHora,
        $AwayOpeningOdds,
        $AwayOpeningDataHora,
        $AwayClosingOdds,
        $AwayClosingDataHora,
        $Url
    )

    $ws.Cells.Item($RowNum, 1).Value2 = $Indice
    $ws.Cells.Item($RowNum, 2).Value2 = "indonesia"
    $ws.Cells.Item($RowNum, 3).Value2 = "liga-1"
    $ws.Cells.Item($RowNum, 4).Value2 = "2023-2024"
    $ws.Cells.Item($RowNum, 5).Value2 = $DataPartida
    $ws.Cells.Item($RowNum, 6).Value2 = $Home
    $ws.Cells.Item($RowNum, 7).Value2 = $HomeGoals
    $ws.Cells.Item($RowNum, 8).Value2 = $Away
    $ws.Cells.Item($RowNum, 9).Value2 = $AwayGoals
    $ws.Cells.Item($RowNum, 10).Value2 = $HomeOpeningOdds
    $ws.Cells.Item($RowNum, 11).Value2 = $HomeOpeningDataHora
    $ws.Cells.Item($RowNum, 12).Value2 = $HomeClosingOdds
    $ws.Cells.Item($RowNum, 13).Value2 = $HomeClosingDataHora
    $ws.Cells.Item($RowNum, 14).Value2 = $DrawOpeningOdds
    $ws.Cells.Item($RowNum, 15).Value2 = $DrawOpeningDataHora
    $ws.Cells.Item($RowNum, 16).Value2 = $DrawClosingOdds
    $ws.Cells.Item($RowNum, 17).Value2 = $DrawClosingDataHora
    $ws.Cells.Item($RowNum, 18).Value2 = $AwayOpeningOdds
    $ws.Cells.Item($RowNum, 19).Value2 = $AwayOpeningDataHora
    $ws.Cells.Item($RowNum, 20).Value2 = $AwayClosingOdds
    $ws.Cells.Item($RowNum, 21).Value2 = $AwayClosingDataHora
    $ws.Cells.Item($RowNum, 22).Value2 = $Url
}

Set-MatchRow 165 164 "Barito Putera" 2 "Persebaya" 0 45239.375 `
    1.93 "07/11/2023 21:12" 1.85 "09/11/2023 08:50" `
    3.49 "07/11/2023 21:12" 3.42 "09/11/2023 08:59" `
    3.47 "07/11/2023 21:12" 4.56 "09/11/2023 08:59" `
    "https://www.betexplorer.com/football/indonesia/liga-1/ps-barito-putera-persebaya/UyhdN4BO/"

Set-MatchRow 166 165 "PSIS Semarang" 4 "Persita" 0 45239.375 `
    1.66 "07/11/2023 21:12" 1.52 "09/11/2023 08:57" `
    3.56 "07/11/2023 21:12" 4.06 "09/11/2023 08:58" `
    4.77 "07/11/2023 21:12" 6.53 "09/11/2023 08:57" `
    "https://www.betexplorer.com/football/indonesia/liga-1/psis-semarang-persita/WtWvnllP/"
